# Update "want to go" counts (column F) on three worksheets to reflect
# newly generated output data (gh-pages rebuild at commit 456a3b4).

$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions)
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F3").Value = 495
$wsExpo.Range("F9").Value = 1552

# Sheet "演出" (performances)
$wsShow = $wb.Worksheets.Item("演出")
$wsShow.Range("F2").Value = 90

# Sheet "全部类型" (all types - combined listing)
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F3").Value = 90
$wsAll.Range("F4").Value = 495
$wsAll.Range("F10").Value = 1552
